# Updated cryptos list on Sat Feb 17 08:40:57 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row
# with the latest scraped values. D-column updates temporarily force a
# text NumberFormat so values like "109.93" are stored as text (matching
# the original inline-string cells) rather than being auto-parsed into
# numbers; the style is reset to Normal afterward so no stray formatting
# is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "51.779.62"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.08%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.780.65"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "357.06"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "109.93"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.56%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.557"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.57%  "
$ws.Range("E8").Value = "  +0.02%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.588"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.59%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "39.88"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.06%  "
$ws.Range("E11").Value = "  +2.16%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0845"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.27%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "19.49"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.04%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.59"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.41%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.223.34"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.07%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.781.05"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.95%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.938"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.05%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "51.721.44"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.42"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.87%  "
$ws.Range("E20").Value = "  -2.33%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.19"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.32%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.0₃0971"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.25%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "70.24"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "269.23"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("E25").Value = "  -2.09%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "26.38"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  +17.22%  "
$ws.Range("E29").Value = "  -0.57%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.25"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.08%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.28"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +7.09%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "51.96"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.22%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "34.82"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.26%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0453"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -8.66%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0841"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.11"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -6.07%  "
$ws.Range("E37").Value = "  +0.09%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "18.68"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.77%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.13"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.90%  "
$ws.Range("E40").Value = "  -3.60%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.56"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("E43").Value = "  -2.04%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "119.97"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -4.72%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "21.68"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -6.45%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.081.34"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.76%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.27"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("E48").Value = "  +1.05%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "5.73"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -4.50%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.938"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.86%  "
$ws.Range("E51").Value = "  +1.64%  "
